$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'37.367.11"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -1.28%  "
$c = $ws.Range("D3")
$c.Value = "'2.050.71"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -1.55%  "
$ws.Range("E4").Value = "  +0.14%  "
$c = $ws.Range("D5")
$c.Value = "'230.95"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.97%  "
$c = $ws.Range("D8")
$c.Value = "'56.99"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -4.01%  "
$ws.Range("E9").Value = "  -2.83%  "
$c = $ws.Range("D10")
$c.Value = "'0.0771"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.32%  "
$c = $ws.Range("D12")
$c.Value = "'2.352.28"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.59%  "
$c = $ws.Range("D13")
$c.Value = "'14.60"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -1.10%  "
$c = $ws.Range("D14")
$c.Value = "'20.64"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -2.73%  "
$ws.Range("E15").Value = "  -2.50%  "
$ws.Range("E16").Value = "  -1.66%  "
$c = $ws.Range("D17")
$c.Value = "'2.042.35"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -1.68%  "
$c = $ws.Range("D18")
$c.Value = "'37.308.27"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.18%  "
$c = $ws.Range("D19")
$c.Value = "'6.09"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.27%  "
$c = $ws.Range("D20")
$c.Value = "'69.64"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -2.66%  "
$ws.Range("E21").Value = "  -3.22%  "
$c = $ws.Range("D22")
$c.Value = "'226.24"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("E24").Value = "  +0.29%  "
$ws.Range("E25").Value = "  -3.69%  "
$c = $ws.Range("D26")
$c.Value = "'9.84"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +7.31%  "
$c = $ws.Range("D27")
$c.Value = "'170.06"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.76%  "
$ws.Range("E28").Value = "  -6.09%  "
$ws.Range("E29").Value = "  -1.60%  "
$ws.Range("E30").Value = "  -5.40%  "
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("E32").Value = "  -4.33%  "
$ws.Range("E33").Value = "  -1.69%  "
$ws.Range("E34").Value = "  -3.95%  "
$ws.Range("E35").Value = "  -1.58%  "
$ws.Range("E36").Value = "  +0.16%  "
$c = $ws.Range("D37")
$c.Value = "'3.27"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -4.94%  "
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("E39").Value = "  -1.97%  "
$c = $ws.Range("D40")
$c.Value = "'0.0224"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +3.22%  "
$c = $ws.Range("D41")
$c.Value = "'98.10"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.91%  "
$c = $ws.Range("D42")
$c.Value = "'0.0953"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.17%  "
$ws.Range("E43").Value = "  +0.16%  "
$c = $ws.Range("D44")
$c.Value = "'1.476.60"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +2.18%  "
$ws.Range("E45").Value = "  +2.12%  "
$c = $ws.Range("D46")
$c.Value = "'16.61"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.59%  "
$c = $ws.Range("D47")
$c.Value = "'3.98"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -4.78%  "
$ws.Range("E48").Value = "  -3.15%  "
$ws.Range("E49").Value = "  -2.10%  "
$ws.Range("E50").Value = "  -2.09%  "
$c = $ws.Range("D51")
$c.Value = "'2.239.08"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -1.55%  "
